# Append three new rows (227, 228, 229) of data to the sheet, matching
# the style/format of the preceding rows (e.g. row 226).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to append: [date serial, B, C, D]
$newRows = @(
    @(44301, 1, 11, 294.4325481798715),
    @(44302, 2, 8,  214.1327623126338),
    @(44303, 1, 6,  160.5995717344754)
)

$startRow = 227

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Clone the formatting of the row above (same layout for every data row)
    # before writing the new values, so the new cells pick up the exact same
    # style (borders, bold, alignment, number format) as the existing ones.
    $ws.Range($ws.Cells.Item($r - 1, 1), $ws.Cells.Item($r - 1, 4)).Copy()
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 4)).PasteSpecial(-4122)

    # Column A: date value.
    $ws.Cells.Item($r, 1).Value = $data[0]

    # Columns B, C, D: plain numeric values.
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
